# Apply the "Succes.xlsx" update described by the commit:
#   - Rename "Coupe du rythme" -> "Coupe Rythme" and
#     "Coupe du challenge" -> "Coupe Challenge"
#   - Add ten new achievement rows (Confrontation / Free For All / Tournoi Score /
#     Tournoi Point / Elimination / Victoire / Victoire assumee / Victoire totale /
#     defi solo / sightread) on Feuil1
#   - Widen column B to fit the new (longer) descriptions
#   - Move the active selection to B37

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Feuil1")

# --- Rename the two pre-existing "online" placeholder entries ---------------
$ws.Range("A3").Value = "Coupe Rythme"
$ws.Range("A4").Value = "Coupe Challenge"

# --- New rows 29-38 -----------------------------------------------------
$newRows = @(
  @{ Row = 29; A = "Coupe Confrontation";    B = "Affrontez des joueurs en ligne";                              C = "Nombre de partie online";        D = 10; E = 50; F = 200; G = 500 },
  @{ Row = 30; A = "Coupe Free For All";     B = "Gagnez des matchs";                                           C = "Nombre de victoire";             D = 5;  E = 25; F = 100; G = 200 },
  @{ Row = 31; A = "Coupe Tournoi Score";    B = "Remportez un tournoi score";                                  C = "Nombre de participant";          D = 3;  E = 5;  F = 7;   G = 8 },
  @{ Row = 32; A = "Coupe Tournoi Point";    B = "Remportez un tournoi point";                                  C = "Nombre de participant";          D = 3;  E = 5;  F = 7;   G = 8 },
  @{ Row = 33; A = "Coupe Elimination";      B = "Remportez une élimination";                                   C = "Nombre de participant";          D = 3;  E = 5;  F = 7;   G = 8 },
  @{ Row = 34; A = "Coupe Victoire";         B = "Sortez victorieux (non FFA)";                                 C = "Nombre de victoire finale";      D = 5;  E = 10; F = 35;  G = 70 },
  @{ Row = 35; A = "Coupe Victoire assumée"; B = "Remportez un FFA à 8";                                        C = "Nombre de victoire";             D = 1;  E = 5;  F = 20;  G = 50 },
  @{ Row = 36; A = "Coupe Victoire totale";  B = "Remportez un tournoi à 8";                                    C = "Nombre de victoire";             D = 2;  E = 5;  F = 15;  G = 30 },
  @{ Row = 37; A = "Coupe défi solo";        B = "Battez un score qui ne vous appartient pas";                  C = "Nombre de score battus";         D = 1;  E = 10; F = 50;  G = 100 },
  @{ Row = 38; A = "Coupe sightread";        B = "Battez un score du premier coup sur une chanson jamais jouée"; C = "Nombre de scores battus";       D = 1;  E = 5;  F = 10;  G = 20 }
)

foreach ($r in $newRows) {
    $ws.Cells.Item($r.Row, 1).Value = $r.A
    $ws.Cells.Item($r.Row, 2).Value = $r.B
    $ws.Cells.Item($r.Row, 3).Value = $r.C
    $ws.Cells.Item($r.Row, 4).Value = $r.D
    $ws.Cells.Item($r.Row, 5).Value = $r.E
    $ws.Cells.Item($r.Row, 6).Value = $r.F
    $ws.Cells.Item($r.Row, 7).Value = $r.G
}

# --- Column B is now much wider to accommodate the longer text --------------
$ws.Columns("B").ColumnWidth = 62.6

# --- Move the selection to B37, matching the author's last edit position ----
$ws.Range("B37").Select()
